# Updated symbol list on Sat Dec 24 19:47:46 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates -- values are stored as text in the sheet, so a
# leading apostrophe is used to force text entry and preserve the exact
# digit/decimal formatting (e.g. trailing zeros) instead of Excel's default
# numeric auto-conversion.
$ws.Range("D2").Value  = "'244.70"
$ws.Range("D3").Value  = "'21.88"
$ws.Range("D5").Value  = "'0.06038"
$ws.Range("D6").Value  = "'3.392"
$ws.Range("D7").Value  = "'0.8147"
$ws.Range("D8").Value  = "'0.9299"
$ws.Range("D10").Value = "'0.07500"
$ws.Range("D12").Value = "'0.03046"
$ws.Range("D13").Value = "'0.09443"
$ws.Range("D14").Value = "'4.016"
$ws.Range("D16").Value = "'0.04814"
$ws.Range("D17").Value = "'0.0005941"
$ws.Range("D18").Value = "'0.005470"
$ws.Range("D19").Value = "'0.004168"
$ws.Range("D20").Value = "'0.0009888"
$ws.Range("D24").Value = "'0.3249"
$ws.Range("D26").Value = "'0.00007001"
$ws.Range("D40").Value = "'0.03998"
$ws.Range("D41").Value = "'0.006411"
$ws.Range("D42").Value = "'0.1078"
$ws.Range("D43").Value = "'0.002720"
$ws.Range("D44").Value = "'0.005910"
$ws.Range("D45").Value = "'0.00005251"
$ws.Range("D47").Value = "'1.000"
$ws.Range("D48").Value = "'0.002333"

# Volume(1h) label (column E) updates
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
